$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 16 (pushes existing rows 16-41 down to 17-42).
$ws.Range("A16").EntireRow.Insert()

$ws.Range("A16").Value = 3
$ws.Range("B16").Value = "Femacal de La Calera"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44874
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 300000000
$ws.Range("G16").Value = "Espárragos"
$ws.Range("H16").Value = "Verde"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 1310
$ws.Range("K16").Value = 1400
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = 1450
$ws.Range("N16").Value = "$/kilo"
$ws.Range("O16").Value = "Provincia de Quillota"
$ws.Range("P16").Value = 1450
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"
